$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 25 (Cambridge 9 Test 1): fill in the Writing score
$ws.Range("I25").Value = 1.1

# Row 27 (new entry: Cambridge 7 Test 2) - practiced Writing
$ws.Range("C27").Value = 45492
$ws.Range("C27").NumberFormat = "[$-409]d\-mmm\-yyyy;@"

$ws.Range("D27").Value = "IELTS7_Test2"

$ws.Range("E27").Value = 30

$ws.Range("F27").Style = "Heading 3"
$ws.Range("F27").NumberFormat = "0.0"
$ws.Range("F27").Formula = "=IFERROR(INDEX(Sheet2!`$F`$5:`$F`$20, MATCH(Table1[[#This Row],[Lis_Mark]], Sheet2!`$D`$5:`$D`$20, 1)),""No Grade"")"

$ws.Range("I27").Value = 1.1

# Move the active selection to match where the user ended up
$ws.Range("C28").Select()
